$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.933.61'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '1.767.93'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '328.75'
$ws.Range('E5').Value = '  +1.35%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4562'
$ws.Range('E7').Value = '  +1.58%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3525'
$ws.Range('E8').Value = '  -1.11%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '42.06'
$ws.Range('E9').Value = '  +1.42%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07386'
$ws.Range('E10').Value = '  -1.07%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.095'
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.70'
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.002'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.192'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('D16').Value = '1.772.00'
$ws.Range('E16').Value = '  +0.90%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '92.60'
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001060'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06444'
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.93'
$ws.Range('E21').Value = '  -1.33%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.771'
$ws.Range('E22').Value = '  +0.65%  '
$ws.Range('D23').Value = '27.962.51'
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.21'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.150'
$ws.Range('E25').Value = '  +2.95%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '161.07'
$ws.Range('E26').Value = '  -2.85%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.15'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').Value = '1.976.70'
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.170'
$ws.Range('E29').Value = '  +3.68%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '124.09'
$ws.Range('E30').Value = '  -1.22%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.076'
$ws.Range('E31').Value = '  -1.71%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09268'
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.663'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.593'
$ws.Range('E34').Value = '  +1.48%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '11.83'
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02276'
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06128'
$ws.Range('E37').Value = '  +1.77%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2086'
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.932'
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6253'
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.182'
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.380'
$ws.Range('E42').Value = '  -0.75%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.824'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5852'
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '122.50'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.935'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.129'
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06818'
$ws.Range('E50').Value = '  -0.92%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '73.32'
$ws.Range('E51').Value = '  +2.37%  '
